$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Label text capitalization fix (added to shared strings first)
$ws.Range("A24").Value = "Backup_DC_Tunnel_Speed"

# Region changed from EMEA to NAM (added to shared strings second)
$ws.Range("B2").Value = "NAM"

# Backup DC tunnel speed value changed from 20 to 40
$ws.Range("B24").Value = 40

# Update selection / view: topLeftCell back to A1 (default), selection now D7
$ws.Range("D7").Select()
